$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Cells.Item(18, 8).Value = 5549.8335
$ws.Cells.Item(18, 10).Value = 8666.333000000001
$ws.Cells.Item(18, 12).Value = 8666.333000000001
$ws.Cells.Item(18, 14).Value = -9234.333000000001
# row 43
$ws.Cells.Item(43, 8).Value = 2448.8572
$ws.Cells.Item(43, 9).Value = 2114.6667
$ws.Cells.Item(43, 10).Value = 2699.5
$ws.Cells.Item(43, 11).Value = 2114.6667
$ws.Cells.Item(43, 12).Value = 2699.5
$ws.Cells.Item(43, 13).Value = -2045.6667
$ws.Cells.Item(43, 14).Value = -2837.5
# row 80
$ws.Cells.Item(80, 8).Value = 1453871.2
$ws.Cells.Item(80, 10).Value = 7239.8335
$ws.Cells.Item(80, 12).Value = 21719.5005
$ws.Cells.Item(80, 14).Value = -23715.5005
# row 82
$ws.Cells.Item(82, 8).Value = 3333
$ws.Cells.Item(82, 9).Value = 3333
$ws.Cells.Item(82, 11).Value = 9999
$ws.Cells.Item(82, 13).Value = -9593
# row 83
$ws.Cells.Item(83, 8).Value = 1453871.2
$ws.Cells.Item(83, 10).Value = 7239.8335
$ws.Cells.Item(83, 12).Value = 65158.5015
$ws.Cells.Item(83, 14).Value = -75142.5015
# row 85
$ws.Cells.Item(85, 8).Value = 3333
$ws.Cells.Item(85, 9).Value = 3333
$ws.Cells.Item(85, 11).Value = 9999
$ws.Cells.Item(85, 13).Value = -8595
# row 86
$ws.Cells.Item(86, 8).Value = 3654
$ws.Cells.Item(86, 10).Value = 3784.8
$ws.Cells.Item(86, 12).Value = 3784.8
$ws.Cells.Item(86, 14).Value = -6030.8
# row 88
$ws.Cells.Item(88, 8).Value = 7198.8
$ws.Cells.Item(88, 10).Value = 8570
$ws.Cells.Item(88, 12).Value = 8570
$ws.Cells.Item(88, 14).Value = -9382
# row 89
$ws.Cells.Item(89, 8).Value = 3654
$ws.Cells.Item(89, 10).Value = 3784.8
$ws.Cells.Item(89, 12).Value = 18924
$ws.Cells.Item(89, 14).Value = -30156
# row 91
$ws.Cells.Item(91, 8).Value = 7198.8
$ws.Cells.Item(91, 10).Value = 8570
$ws.Cells.Item(91, 12).Value = 8570
$ws.Cells.Item(91, 14).Value = -11378
# row 107
$ws.Cells.Item(107, 8).Value = 827.26086
$ws.Cells.Item(107, 9).Value = 936.4
$ws.Cells.Item(107, 11).Value = 936.4
$ws.Cells.Item(107, 13).Value = 983.6
# row 113
$ws.Cells.Item(113, 8).Value = 6785.75
$ws.Cells.Item(113, 10).Value = 6464.3335
$ws.Cells.Item(113, 12).Value = 6464.3335
$ws.Cells.Item(113, 14).Value = -12972.3335
# row 132
$ws.Cells.Item(132, 8).Value = 7938601.5
$ws.Cells.Item(132, 10).Value = 1194
$ws.Cells.Item(132, 12).Value = 3582
$ws.Cells.Item(132, 14).Value = -8642
# row 138
$ws.Cells.Item(138, 8).Value = 1419897.1
$ws.Cells.Item(138, 9).Value = 2355.0588
$ws.Cells.Item(138, 10).Value = 2054060.8
$ws.Cells.Item(138, 11).Value = 7065.176399999999
$ws.Cells.Item(138, 12).Value = 6162182.4
$ws.Cells.Item(138, 13).Value = -1925.176399999999
$ws.Cells.Item(138, 14).Value = -6172462.4
# row 141
$ws.Cells.Item(141, 8).Value = 2126.1177
$ws.Cells.Item(141, 9).Value = 2126.1177
$ws.Cells.Item(141, 11).Value = 6378.353099999999
$ws.Cells.Item(141, 13).Value = -1198.353099999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Cells.Item(32, 8).Value = 18646.367
$ws.Cells.Item(32, 9).Value = 20836.547
$ws.Cells.Item(32, 11).Value = 20836.547
$ws.Cells.Item(32, 13).Value = -20549.547
# row 34
$ws.Cells.Item(34, 8).Value = 59000
# row 61
$ws.Cells.Item(61, 8).Value = 5954.4614
$ws.Cells.Item(61, 9).Value = 3406.7878
$ws.Cells.Item(61, 11).Value = 3406.7878
$ws.Cells.Item(61, 13).Value = -3194.7878
# row 74
$ws.Cells.Item(74, 8).Value = 3395.875
$ws.Cells.Item(74, 9).Value = 1337.125
$ws.Cells.Item(74, 10).Value = 8542.75
$ws.Cells.Item(74, 11).Value = 1337.125
$ws.Cells.Item(74, 12).Value = 8542.75
$ws.Cells.Item(74, 13).Value = -463.125
$ws.Cells.Item(74, 14).Value = -10290.75
# row 77
$ws.Cells.Item(77, 8).Value = 3395.875
$ws.Cells.Item(77, 9).Value = 1337.125
$ws.Cells.Item(77, 10).Value = 8542.75
$ws.Cells.Item(77, 11).Value = 6685.625
$ws.Cells.Item(77, 12).Value = 42713.75
$ws.Cells.Item(77, 13).Value = -2317.625
$ws.Cells.Item(77, 14).Value = -51449.75
# row 97
$ws.Cells.Item(97, 8).Value = 1007.2381
$ws.Cells.Item(97, 9).Value = 836
$ws.Cells.Item(97, 10).Value = 1349.7142
$ws.Cells.Item(97, 11).Value = 836
$ws.Cells.Item(97, 12).Value = 1349.7142
$ws.Cells.Item(97, 13).Value = -340
$ws.Cells.Item(97, 14).Value = -2341.7142
# row 132
$ws.Cells.Item(132, 8).Value = 1561.3636
$ws.Cells.Item(132, 9).Value = 1311.711
$ws.Cells.Item(132, 11).Value = 3935.133
$ws.Cells.Item(132, 13).Value = -1405.133
# row 136
$ws.Cells.Item(136, 8).Value = 5954.4614
$ws.Cells.Item(136, 9).Value = 3406.7878
$ws.Cells.Item(136, 11).Value = 10220.3634
$ws.Cells.Item(136, 13).Value = -7670.3634

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Cells.Item(94, 8).Value = 6250893
$ws.Cells.Item(94, 9).Value = 1063.091
$ws.Cells.Item(94, 11).Value = 1063.091
$ws.Cells.Item(94, 13).Value = -612.0909999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Cells.Item(7, 8).Value = 210.11539
$ws.Cells.Item(7, 9).Value = 197.94737
$ws.Cells.Item(7, 10).Value = 243.14285
$ws.Cells.Item(7, 11).Value = 197.94737
$ws.Cells.Item(7, 12).Value = 243.14285
$ws.Cells.Item(7, 13).Value = -84.94737000000001
$ws.Cells.Item(7, 14).Value = -469.14285
# row 22
$ws.Cells.Item(22, 8).Value = 358.9
$ws.Cells.Item(22, 9).Value = 256
$ws.Cells.Item(22, 10).Value = 513.25
$ws.Cells.Item(22, 11).Value = 256
$ws.Cells.Item(22, 12).Value = 513.25
$ws.Cells.Item(22, 13).Value = 94
$ws.Cells.Item(22, 14).Value = -1213.25
# row 31
$ws.Cells.Item(31, 8).Value = 1925545.4
$ws.Cells.Item(31, 9).Value = 4764350
$ws.Cells.Item(31, 10).Value = 2484.3225
$ws.Cells.Item(31, 11).Value = 4764350
$ws.Cells.Item(31, 12).Value = 2484.3225
$ws.Cells.Item(31, 13).Value = -4764055
$ws.Cells.Item(31, 14).Value = -3074.3225
# row 34
$ws.Cells.Item(34, 8).Value = 1925545.4
$ws.Cells.Item(34, 9).Value = 4764350
$ws.Cells.Item(34, 10).Value = 2484.3225
$ws.Cells.Item(34, 11).Value = 4764350
$ws.Cells.Item(34, 12).Value = 2484.3225
$ws.Cells.Item(34, 13).Value = -4764148
$ws.Cells.Item(34, 14).Value = -2888.3225
# row 58
$ws.Cells.Item(58, 8).Value = 1660.65
$ws.Cells.Item(58, 9).Value = 1016.8
$ws.Cells.Item(58, 10).Value = 2304.5
$ws.Cells.Item(58, 11).Value = 1016.8
$ws.Cells.Item(58, 12).Value = 2304.5
$ws.Cells.Item(58, 13).Value = -813.8
$ws.Cells.Item(58, 14).Value = -2710.5
# row 132
$ws.Cells.Item(132, 8).Value = 1929.6945
$ws.Cells.Item(132, 9).Value = 1471.2727
$ws.Cells.Item(132, 11).Value = 4413.8181
$ws.Cells.Item(132, 13).Value = -1883.8181
# row 136
$ws.Cells.Item(136, 8).Value = 1660.65
$ws.Cells.Item(136, 9).Value = 1016.8
$ws.Cells.Item(136, 10).Value = 2304.5
$ws.Cells.Item(136, 11).Value = 3050.4
$ws.Cells.Item(136, 12).Value = 6913.5
$ws.Cells.Item(136, 13).Value = -500.3999999999996
$ws.Cells.Item(136, 14).Value = -12013.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 8
$ws.Cells.Item(8, 8).Value = 650
$ws.Cells.Item(8, 9).Value = 650
$ws.Cells.Item(8, 11).Value = 1950
$ws.Cells.Item(8, 13).Value = -1811

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Cells.Item(70, 8).Value = 3999.6667
$ws.Cells.Item(70, 9).Value = 3999.6667
$ws.Cells.Item(70, 11).Value = 3999.6667
$ws.Cells.Item(70, 13).Value = -3729.6667
# row 73
$ws.Cells.Item(73, 8).Value = 3999.6667
$ws.Cells.Item(73, 9).Value = 3999.6667
$ws.Cells.Item(73, 11).Value = 3999.6667
$ws.Cells.Item(73, 13).Value = -3063.6667
# row 132
$ws.Cells.Item(132, 8).Value = 2496.7932
$ws.Cells.Item(132, 9).Value = 2326.4546
$ws.Cells.Item(132, 10).Value = 3032.1428
$ws.Cells.Item(132, 11).Value = 6979.3638
$ws.Cells.Item(132, 12).Value = 9096.428400000001
$ws.Cells.Item(132, 13).Value = -4449.3638
$ws.Cells.Item(132, 14).Value = -14156.4284

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Cells.Item(7, 8).Value = 4725
$ws.Cells.Item(7, 9).Value = 4000
$ws.Cells.Item(7, 10).Value = 4966.6665
$ws.Cells.Item(7, 11).Value = 4000
$ws.Cells.Item(7, 12).Value = 4966.6665
$ws.Cells.Item(7, 13).Value = -3888
$ws.Cells.Item(7, 14).Value = -5190.6665
# row 122
$ws.Cells.Item(122, 8).Value = 3500
$ws.Cells.Item(122, 9).Value = 4000
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 13).Value = -9550
# row 126
$ws.Cells.Item(126, 8).Value = 4725
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 10).Value = 4966.6665
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 14899.9995
$ws.Cells.Item(126, 13).Value = -9530
$ws.Cells.Item(126, 14).Value = -19839.9995
# row 132
$ws.Cells.Item(132, 8).Value = 5524.484
$ws.Cells.Item(132, 9).Value = 5634.7144
$ws.Cells.Item(132, 10).Value = 5433.706
$ws.Cells.Item(132, 11).Value = 16904.1432
$ws.Cells.Item(132, 12).Value = 16301.118
$ws.Cells.Item(132, 13).Value = -14374.1432
$ws.Cells.Item(132, 14).Value = -21361.118
# row 136
$ws.Cells.Item(136, 8).Value = 4103.6
$ws.Cells.Item(136, 9).Value = 1869.5
$ws.Cells.Item(136, 11).Value = 5608.5
$ws.Cells.Item(136, 13).Value = -3058.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Cells.Item(113, 8).Value = 713.1795
$ws.Cells.Item(113, 9).Value = 758.1875
$ws.Cells.Item(113, 10).Value = 507.42856
$ws.Cells.Item(113, 11).Value = 2274.5625
$ws.Cells.Item(113, 12).Value = 1522.28568
$ws.Cells.Item(113, 13).Value = -104.5625
$ws.Cells.Item(113, 14).Value = -5862.28568
# row 132
$ws.Cells.Item(132, 8).Value = 19613.125
$ws.Cells.Item(132, 9).Value = 24447.973
$ws.Cells.Item(132, 11).Value = 73343.91900000001
$ws.Cells.Item(132, 13).Value = -70813.91900000001
# row 136
$ws.Cells.Item(136, 8).Value = 8361.333000000001
$ws.Cells.Item(136, 9).Value = 10903.292
$ws.Cells.Item(136, 10).Value = 4972.0557
$ws.Cells.Item(136, 11).Value = 32709.876
$ws.Cells.Item(136, 12).Value = 14916.1671
$ws.Cells.Item(136, 13).Value = -30159.876
$ws.Cells.Item(136, 14).Value = -20016.1671
